$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 1249.5
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1249.5
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1249.5
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -1589.5
# Row 17
$ws.Range("H17").Value = 2558.9312
$ws.Range("J17").Value = 2558.9312
$ws.Range("L17").Value = 7676.7936
$ws.Range("N17").Value = -8012.7936
# Row 18
$ws.Range("H18").Value = 2086.375
$ws.Range("I18").Value = 1598.7142
$ws.Range("K18").Value = 1598.7142
$ws.Range("M18").Value = -1314.7142
# Row 112
$ws.Range("H112").Value = 1845.619
$ws.Range("J112").Value = 1845.619
$ws.Range("L112").Value = 5536.857
$ws.Range("N112").Value = -7752.857
# Row 127
$ws.Range("H127").Value = 931
$ws.Range("I127").Value = 1004.8333
$ws.Range("J127").Value = 783.3333
$ws.Range("K127").Value = 3014.4999
$ws.Range("L127").Value = 2349.9999
$ws.Range("M127").Value = 1945.5001
$ws.Range("N127").Value = -12269.9999
# Row 132
$ws.Range("H132").Value = 783.8823
$ws.Range("I132").Value = 801
$ws.Range("K132").Value = 2403
$ws.Range("M132").Value = 127
# Row 137
$ws.Range("H137").Value = 6950397
$ws.Range("I137").Value = 11366594
$ws.Range("J137").Value = 10659.357
$ws.Range("K137").Value = 34099782
$ws.Range("L137").Value = 31978.071
$ws.Range("M137").Value = -34097232
$ws.Range("N137").Value = -37078.071
# Row 138
$ws.Range("H138").Value = 2996.814
$ws.Range("I138").Value = 2186
$ws.Range("J138").Value = 3770.7727
$ws.Range("K138").Value = 6558
$ws.Range("L138").Value = 11312.3181
$ws.Range("M138").Value = -1418
$ws.Range("N138").Value = -21592.3181
# Row 141
$ws.Range("H141").Value = 1979.8182
$ws.Range("I141").Value = 1953.2222
$ws.Range("K141").Value = 5859.6666
$ws.Range("M141").Value = -679.6665999999996

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10066366
$ws.Range("I32").Value = 11841202
$ws.Range("K32").Value = 11841202
$ws.Range("M32").Value = -11840915
# Row 45
$ws.Range("H45").Value = 27116.889
$ws.Range("I45").Value = 28830.895
$ws.Range("K45").Value = 28830.895
$ws.Range("M45").Value = -28453.895
# Row 74
$ws.Range("H74").Value = 336203.84
$ws.Range("I74").Value = 436692.78
$ws.Range("J74").Value = 6025.857
$ws.Range("K74").Value = 436692.78
$ws.Range("L74").Value = 6025.857
$ws.Range("M74").Value = -435818.78
$ws.Range("N74").Value = -7773.857
# Row 77
$ws.Range("H77").Value = 336203.84
$ws.Range("I77").Value = 436692.78
$ws.Range("J77").Value = 6025.857
$ws.Range("K77").Value = 2183463.9
$ws.Range("L77").Value = 30129.285
$ws.Range("M77").Value = -2179095.9
$ws.Range("N77").Value = -38865.285
# Row 92
$ws.Range("H92").Value = 39270
$ws.Range("J92").Value = 39270
$ws.Range("L92").Value = 39270
$ws.Range("N92").Value = -44262
# Row 121
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
# Row 132
$ws.Range("H132").Value = 4000.537
$ws.Range("I132").Value = 2788.558
$ws.Range("K132").Value = 8365.673999999999
$ws.Range("M132").Value = -5835.673999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 3898.5
$ws.Range("J22").Value = 5499.75
$ws.Range("L22").Value = 5499.75
$ws.Range("N22").Value = -5845.75
# Row 134
$ws.Range("H134").Value = 3348.3396
$ws.Range("I134").Value = 2399.3635
$ws.Range("K134").Value = 7198.0905
$ws.Range("M134").Value = -4663.0905

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 733
$ws.Range("I7").Value = 734.40625
$ws.Range("K7").Value = 734.40625
$ws.Range("M7").Value = -621.40625
# Row 19
$ws.Range("H19").Value = 985.6857
$ws.Range("I19").Value = 920.0417
$ws.Range("K19").Value = 920.0417
$ws.Range("M19").Value = -750.0417
# Row 22
$ws.Range("H22").Value = 1569.7727
$ws.Range("I22").Value = 701.36365
$ws.Range("J22").Value = 2438.182
$ws.Range("K22").Value = 701.36365
$ws.Range("L22").Value = 2438.182
$ws.Range("M22").Value = -351.36365
$ws.Range("N22").Value = -3138.182
# Row 24
$ws.Range("H24").Value = 985.6857
$ws.Range("I24").Value = 920.0417
$ws.Range("K24").Value = 920.0417
$ws.Range("M24").Value = -750.0417
# Row 41
$ws.Range("H41").Value = 7618.1
$ws.Range("J41").Value = 13572.25
$ws.Range("L41").Value = 13572.25
$ws.Range("N41").Value = -14428.25
# Row 58
$ws.Range("H58").Value = 7983.143
$ws.Range("I58").Value = 5187.857
$ws.Range("K58").Value = 5187.857
$ws.Range("M58").Value = -4984.857
# Row 132
$ws.Range("H132").Value = 3439.9092
$ws.Range("I132").Value = 2308.64
$ws.Range("K132").Value = 6925.92
$ws.Range("M132").Value = -4395.92
# Row 136
$ws.Range("H136").Value = 7983.143
$ws.Range("I136").Value = 5187.857
$ws.Range("K136").Value = 15563.571
$ws.Range("M136").Value = -13013.571

$ws = $wb.Worksheets.Item("CUL")
# Row 26
$ws.Range("H26").Value = 3668.257
$ws.Range("I26").Value = 569
$ws.Range("J26").Value = 12621.667
$ws.Range("K26").Value = 1707
$ws.Range("L26").Value = 37865.001
$ws.Range("M26").Value = -1419
$ws.Range("N26").Value = -38441.001
# Row 80
$ws.Range("H80").Value = 3666.3333
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3666.3333
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 10998.9999
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -12870.9999
# Row 83
$ws.Range("H83").Value = 3666.3333
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3666.3333
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 32996.9997
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -42356.9997

$ws = $wb.Worksheets.Item("GSM")
# Row 14
$ws.Range("H14").Value = 1001354.6
$ws.Range("I14").Value = 1334933.1
$ws.Range("J14").Value = 619
$ws.Range("K14").Value = 1334933.1
$ws.Range("L14").Value = 619
$ws.Range("M14").Value = -1334765.1
$ws.Range("N14").Value = -955
# Row 122
$ws.Range("H122").Value = 3861.0667
$ws.Range("I122").Value = 2452.3
$ws.Range("J122").Value = 6678.6
$ws.Range("K122").Value = 7356.900000000001
$ws.Range("L122").Value = 20035.8
$ws.Range("M122").Value = -4906.900000000001
$ws.Range("N122").Value = -24935.8
# Row 132
$ws.Range("H132").Value = 4737.3335
$ws.Range("I132").Value = 3060.1765
$ws.Range("K132").Value = 9180.529500000001
$ws.Range("M132").Value = -6650.529500000001

$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 18199.572
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
# Row 14
$ws.Range("H14").Value = 7184.7
$ws.Range("I14").Value = 3064
$ws.Range("K14").Value = 3064
$ws.Range("M14").Value = -2892
# Row 15
$ws.Range("H15").Value = 18199.572
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
# Row 22
$ws.Range("H22").Value = 3355.1943
$ws.Range("I22").Value = 1841.6875
$ws.Range("K22").Value = 1841.6875
$ws.Range("M22").Value = -1546.6875
# Row 27
$ws.Range("H27").Value = 3355.1943
$ws.Range("I27").Value = 1841.6875
$ws.Range("K27").Value = 1841.6875
$ws.Range("M27").Value = -1734.6875
# Row 43
$ws.Range("H43").Value = 8666.333000000001
$ws.Range("I43").Value = 8000
$ws.Range("K43").Value = 8000
$ws.Range("M43").Value = -7807
# Row 104
$ws.Range("H104").Value = 12170.429
$ws.Range("J104").Value = 12170.429
$ws.Range("L104").Value = 12170.429
$ws.Range("N104").Value = -19158.429
# Row 122
$ws.Range("H122").Value = 25003430
$ws.Range("I122").Value = 31253020
$ws.Range("J122").Value = 5070.5
$ws.Range("K122").Value = 93759060
$ws.Range("L122").Value = 15211.5
$ws.Range("M122").Value = -93756610
$ws.Range("N122").Value = -20111.5
# Row 132
$ws.Range("H132").Value = 3332.9143
$ws.Range("I132").Value = 1550.3462
$ws.Range("K132").Value = 4651.0386
$ws.Range("M132").Value = -2121.0386
# Row 136
$ws.Range("H136").Value = 5175.05
$ws.Range("I136").Value = 5571.7144
$ws.Range("J136").Value = 4961.4614
$ws.Range("K136").Value = 16715.1432
$ws.Range("L136").Value = 14884.3842
$ws.Range("M136").Value = -14165.1432
$ws.Range("N136").Value = -19984.3842

$ws = $wb.Worksheets.Item("WVR")
# Row 9
$ws.Range("H9").Value = 5500
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 5500
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 5500
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -5780
# Row 14
$ws.Range("H14").Value = 5500
$ws.Range("I14").Value = 5000
$ws.Range("K14").Value = 5000
$ws.Range("M14").Value = -4832
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
# Row 43
$ws.Range("H43").Value = 77500
$ws.Range("I43").Value = 70000
$ws.Range("K43").Value = 70000
$ws.Range("M43").Value = -69851
# Row 104
$ws.Range("H104").Value = 22913
$ws.Range("J104").Value = 22913
$ws.Range("L104").Value = 22913
$ws.Range("N104").Value = -29901
# Row 132
$ws.Range("H132").Value = 4522.3022
$ws.Range("I132").Value = 2872.1428
$ws.Range("J132").Value = 7602.6
$ws.Range("K132").Value = 8616.428400000001
$ws.Range("L132").Value = 22807.8
$ws.Range("M132").Value = -6086.428400000001
$ws.Range("N132").Value = -27867.8
